# end of day backup
$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: "HotelCard" cell - merge the spell-check-wrapped "HotelCard"
# run and the trailing-space run into a single run "HotelCard ", and
# drop the now-stale proofErr markers.
# -----------------------------------------------------------------------
$table = $d.Tables.Item(1)
$hcCell = $table.Cell(3, 1)
$hcPara = $hcCell.Range.Paragraphs.Item(1)

# Append a brand-new (proofErr-free) paragraph with the desired text,
# then remove the original paragraph - leaves exactly one clean <w:p>.
$hcPara.Range.InsertAfter([char]13 + "HotelCard ")

$table = $d.Tables.Item(1)
$hcCell = $table.Cell(3, 1)
$hcOldPara = $hcCell.Range.Paragraphs.Item(1)
$hcOldPara.Range.Delete()

# -----------------------------------------------------------------------
# Change 2: "Not Started" -> "In Progress", plus a new follow-up note
# paragraph ("Requires some fiddling of date params") in the same cell,
# keeping the InProg paragraph style.
# -----------------------------------------------------------------------
$table = $d.Tables.Item(1)
$lastRow = $table.Rows.Count
$statusCell = $table.Cell($lastRow, 3)
$statusPara = $statusCell.Range.Paragraphs.Item(1)
$statusPara.Range.Text = "In Progress"

$table = $d.Tables.Item(1)
$statusCell = $table.Cell($lastRow, 3)
$statusPara = $statusCell.Range.Paragraphs.Item(1)
$statusPara.Range.InsertAfter([char]13 + "Requires some fiddling of date params")
